$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.588.38"
$ws.Range("E2").Value = "  +0.90%  "

# Row 3
$ws.Range("D3").Value = "2.429.88"
$ws.Range("E3").Value = "  +0.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.61"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.59%  "

# Row 6
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.38"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +1.93%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +0.93%  "

# Row 10
$ws.Range("E10").Value = "  +0.32%  "

# Row 11
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +1.43%  "

# Row 12
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +1.54%  "

# Row 13
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.81"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +4.74%  "

# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +3.90%  "

# Row 15
$ws.Range("D15").Value = "2.868.28"
$ws.Range("E15").Value = "  +0.70%  "

# Row 16
$ws.Range("D16").Value = "62.466.12"
$ws.Range("E16").Value = "  +0.91%  "

# Row 17
$ws.Range("D17").Value = "2.430.13"
$ws.Range("E17").Value = "  +0.90%  "

# Row 18
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.23"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.21%  "

# Row 19
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.99"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +2.36%  "

# Row 20
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.03"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.74%  "

# Row 21
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.77%  "

# Row 22
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.18"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.82%  "

# Row 24
$ws.Range("E24").Value = "  +4.98%  "

# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "599.45"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +6.23%  "

# Row 26
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.59"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.91%  "

# Row 27
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +7.80%  "

# Row 28
$ws.Range("D28").Value = "2.548.54"
$ws.Range("E28").Value = "  +0.90%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.44"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +3.01%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +3.31%  "

# Row 32
$ws.Range("E32").Value = "  -3.13%  "

# Row 33
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.43%  "

# Row 34
$ws.Range("E34").Value = "  -1.18%  "

# Row 35
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.86"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +1.73%  "

# Row 36
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.382"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +0.50%  "

# Row 38
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.74"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.05%  "

# Row 39
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.62%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +0.95%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.11"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -3.64%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +10.01%  "

# Row 44
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.68"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.06%  "

# Row 45
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.68"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.85%  "

# Row 46
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.57"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +3.43%  "

# Row 48
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.602"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.70%  "

# Row 49
$ws.Range("E49").Value = "  +2.53%  "

# Row 50
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0920"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.02%  "

# Row 51
$ws.Range("E51").Value = "  +3.99%  "
